$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8268491625785828
$ws.Range("B1").Value = 1.284665465354919
$ws.Range("C1").Value = 4.593391418457031
$ws.Range("D1").Value = 4.051531314849854
$ws.Range("E1").Value = 0.5496358275413513
